$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# N2 already carries the "Hyperlink" cell style (xf index 2) from the original workbook. Grab a
# copy of that formatting now, before we touch its contents, so we can re-stamp it onto every
# email cell later without Hyperlinks.Add's side effect of registering a brand-new (duplicate)
# style record for each cell it touches.
$ws.Cells.Item(2, 14).Copy() | Out-Null
$hyperlinkFormatRow = 2
$hyperlinkFormatCol = 14

# Remove the existing hyperlink (old N2 -> helloguba@gmail.com) before rewriting cells so we don't
# carry forward stale relationships.
for ($i = $ws.Hyperlinks.Count; $i -ge 1; $i--) {
    $ws.Hyperlinks.Item($i).Delete()
}

# --- Header row (row 1): insert "Reference Code 2" at J, shift State/Zip/Country left-stacking to K/L/M,
#     keep Email/Phone where they are, and move "Reference Code" out to the new P column. ---
$ws.Cells.Item(1, 10).Value = "Reference Code 2"   # J1
$ws.Cells.Item(1, 11).Value = "Donor State"        # K1
$ws.Cells.Item(1, 12).Value = "Donor ZIP"          # L1
$ws.Cells.Item(1, 13).Value = "Donor Country"      # M1
$ws.Cells.Item(1, 14).Value = "Donor Email"        # N1 (unchanged)
$ws.Cells.Item(1, 15).Value = "Donor Phone"        # O1 (unchanged)
$ws.Cells.Item(1, 16).Value = "Reference Code"     # P1 (new)

# Clear out the old Phone column content for the data rows -- in the new layout there is no
# phone number recorded for any of rows 2-5.
$ws.Range("O2:O5").Clear()

# --- Row 2 ---
$ws.Cells.Item(2, 1).Value  = "AB10958920755555"
$ws.Cells.Item(2, 2).Value  = 44344.184224537035
$ws.Cells.Item(2, 3).Value  = 20.21
$ws.Cells.Item(2, 4).Value  = "unlimited"
$ws.Cells.Item(2, 5).Value  = 2
$ws.Cells.Item(2, 6).Value  = "Sarah"
$ws.Cells.Item(2, 7).Value  = "Kidd"
$ws.Cells.Item(2, 8).Value  = "TEST Rd"
$ws.Cells.Item(2, 9).Value  = "McKenzie"
$ws.Cells.Item(2, 10).Value = "TEST 1"
$ws.Cells.Item(2, 11).Value = "TN"
$ws.Cells.Item(2, 12).Value = 38201
$ws.Cells.Item(2, 13).Value = "United States"
$ws.Cells.Item(2, 14).Value = "newdonor@bethelu.edu"
$ws.Cells.Item(2, 16).Value = "TESTbackontrack.fr.042821"

# --- Row 3 ---
$ws.Cells.Item(3, 1).Value  = "GB195892342"
$ws.Cells.Item(3, 2).Value  = 44344.184224537035
$ws.Cells.Item(3, 3).Value  = 25
$ws.Cells.Item(3, 4).Value  = "unlimited"
$ws.Cells.Item(3, 5).Value  = 2
$ws.Cells.Item(3, 6).Value  = "Martha"
$ws.Cells.Item(3, 7).Value  = "Lemert"
$ws.Cells.Item(3, 8).Value  = "TEST Rd"
$ws.Cells.Item(3, 9).Value  = "Fort Wayne"
$ws.Cells.Item(3, 10).Value = "TEST 2"
$ws.Cells.Item(3, 11).Value = "IN"
$ws.Cells.Item(3, 12).Value = 46845
$ws.Cells.Item(3, 13).Value = "United States"
$ws.Cells.Item(3, 14).Value = "lemertmartha@gmail.com"
$ws.Cells.Item(3, 16).Value = "backontrack.fr.042821"

# --- Row 4 ---
$ws.Cells.Item(4, 1).Value  = "KK195892342"
$ws.Cells.Item(4, 2).Value  = 44344.184224537035
$ws.Cells.Item(4, 3).Value  = 25
$ws.Cells.Item(4, 4).Value  = "unlimited"
$ws.Cells.Item(4, 5).Value  = 2
$ws.Cells.Item(4, 6).Value  = "Martha"
$ws.Cells.Item(4, 7).Value  = "Lemert"
$ws.Cells.Item(4, 8).Value  = "TEST Rd"
$ws.Cells.Item(4, 9).Value  = "Fort Wayne"
$ws.Cells.Item(4, 10).Value = "TEST 3"
$ws.Cells.Item(4, 11).Value = "IN"
$ws.Cells.Item(4, 12).Value = 46845
$ws.Cells.Item(4, 13).Value = "United States"
$ws.Cells.Item(4, 14).Value = "newdonor@gmail.com"
$ws.Cells.Item(4, 16).Value = "RedistrictingEOM.FR.05.27.21."

# --- Row 5 ---
$ws.Cells.Item(5, 1).Value  = "AB10958920755555"
$ws.Cells.Item(5, 2).Value  = 44344.184224537035
$ws.Cells.Item(5, 3).Value  = 20.21
$ws.Cells.Item(5, 4).Value  = "unlimited"
$ws.Cells.Item(5, 5).Value  = 2
$ws.Cells.Item(5, 6).Value  = "Sarah"
$ws.Cells.Item(5, 7).Value  = "Kidd"
$ws.Cells.Item(5, 8).Value  = "TEST Rd"
$ws.Cells.Item(5, 9).Value  = "McKenzie"
$ws.Cells.Item(5, 10).Value = "TEST 1"
$ws.Cells.Item(5, 11).Value = "TN"
$ws.Cells.Item(5, 12).Value = 38201
$ws.Cells.Item(5, 13).Value = "United States"
$ws.Cells.Item(5, 14).Value = "newdonor@bethelu.edu"
$ws.Cells.Item(5, 16).Value = "TESTbackontrack.fr.042821"

# --- Hyperlinks on the Donor Email column (N4, N5, N2 in that creation order to mirror r:id ordering) ---
$ws.Cells.Item(4, 14).Style = "Hyperlink"
$ws.Cells.Item(5, 14).Style = "Hyperlink"
$ws.Cells.Item(2, 14).Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Cells.Item(4, 14), "mailto:newdonor@gmail.com", "", "", "newdonor@gmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(5, 14), "mailto:newdonor@bethelu.edu", "", "", "newdonor@bethelu.edu") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(2, 14), "mailto:newdonor@bethelu.edu", "", "", "newdonor@bethelu.edu") | Out-Null

# --- Sheet view: drop the frozen top-left cell offset, select B14 ---
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B14").Select()

# --- Sort state range now starts at row 3 (header + first data row excluded from the remembered sort range) ---
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B3:B114")) | Out-Null
$ws.Sort.SetRange($ws.Range("A3:O114"))
$ws.Sort.Header = 0
